# Applies the "gh-pages output generated at 456a3b4" update to
# 南宁-漫展信息.xlsx:
#   - sheet "展览"   (Worksheets(1)): update 想去人数(F) counters, and
#     insert a new row for "南宁·熊喵M动漫嘉年华【免费】" before the
#     "第二届北极光动漫展" row.
#   - sheet "演出"   (Worksheets(2)): update 想去人数(F) counters.
#   - sheet "全部类型" (Worksheets(4)): update 想去人数(F) counters, and
#     insert the same new "熊喵M" row in the matching spot.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $text) {
    # Force the cell to stay a text cell (the source data keeps date-like
    # strings such as "2024-08-17" as literal text, not Excel date serials).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

function Style-IndexCell($ws, $addr) {
    # Reproduce the workbook's "index column" style (bold, thin box border,
    # centered/top aligned) used by every cell in column A.
    $r = $ws.Range($addr)
    $r.Font.Bold = $true
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
    $r.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("G2").Value = 60
$ws1.Range("F3").Value = 654
$ws1.Range("F4").Value = 0
$ws1.Range("F5").Value = 5066
$ws1.Range("F7").Value = 9656
$ws1.Range("F8").Value = 252
$ws1.Range("F10").Value = 93

# Insert a brand-new row 11 ("南宁·熊喵M动漫嘉年华【免费】"); this pushes the
# existing rows 11 ("第二届北极光动漫展") and 12 ("万圣漫控嘉年华10") down to
# rows 12 and 13, carrying their own column-A index values along with them.
$ws1.Rows.Item(11).Insert()

$ws1.Range("A11").Value = 10
Style-IndexCell $ws1 "A11"
Set-TextCell $ws1 "B11" "2024-08-17"
$ws1.Range("C11").Value = "南宁·熊喵M动漫嘉年华【免费】"
$ws1.Range("D11").Value = "港航上尧码头(江北大道南100米) 水明漾艺术中心"
$ws1.Range("E11").Value = "2024.08.17 10:00-08.18 17:00"
$ws1.Range("F11").Value = 0
$ws1.Range("G11").Value = 29.9
$ws1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=89145"
$ws1.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202407/ndmB7MOh1720344131003.jpeg"

# ---------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F4").Value = 10
$ws2.Range("F6").Value = 3

# ---------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 0
$ws4.Range("G2").Value = 60
$ws4.Range("F4").Value = 354
$ws4.Range("F7").Value = 5066
$ws4.Range("F8").Value = 533
$ws4.Range("F10").Value = 9656
$ws4.Range("F12").Value = 535

# Insert the new "熊喵M" row before row 16 ("第二届北极光动漫展"); this
# pushes rows 16-18 down to 17-19.
$ws4.Rows.Item(16).Insert()

$ws4.Range("A16").Value = 15
Style-IndexCell $ws4 "A16"
Set-TextCell $ws4 "B16" "2024-08-17"
$ws4.Range("C16").Value = "南宁·熊喵M动漫嘉年华【免费】"
$ws4.Range("D16").Value = "港航上尧码头(江北大道南100米) 水明漾艺术中心"
$ws4.Range("E16").Value = "2024.08.17 10:00-08.18 17:00"
$ws4.Range("F16").Value = 11
$ws4.Range("G16").Value = 29.9
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=89145"
$ws4.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202407/ndmB7MOh1720344131003.jpeg"

# Row 17 (shifted from the old row 16, "第二届北极光动漫展"): only its
# 想去人数(F) count changed.
$ws4.Range("F17").Value = 706

Write-Output "done"
